# Apply updated odds values per the target diff.
# Each Range().Value assignment sets a single cell to its new numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: 31 cell(s) updated
$ws.Range("G9").Value = 3.25
$ws.Range("I9").Value = 2.22
$ws.Range("J9").Value = 3.85
$ws.Range("L9").Value = 2.77
$ws.Range("Q9").Value = 1.91
$ws.Range("R9").Value = 1.7
$ws.Range("X9").Value = 17
$ws.Range("Y9").Value = 11.25
$ws.Range("Z9").Value = 45
$ws.Range("AA9").Value = 32
$ws.Range("AB9").Value = 37
$ws.Range("AG9").Value = 450
$ws.Range("AH9").Value = 7.8
$ws.Range("AI9").Value = 11.25
$ws.Range("AJ9").Value = 8.5
$ws.Range("AK9").Value = 23
$ws.Range("AL9").Value = 17.5
$ws.Range("AM9").Value = 26
$ws.Range("AN9").Value = 5.2
$ws.Range("AO9").Value = 19
$ws.Range("AP9").Value = 25
$ws.Range("AQ9").Value = 100
$ws.Range("AR9").Value = 150
$ws.Range("AS9").Value = 350
$ws.Range("AT9").Value = 2.45
$ws.Range("AW9").Value = 4.1
$ws.Range("AX9").Value = 11.5
$ws.Range("AY9").Value = 18.5
$ws.Range("AZ9").Value = 45
$ws.Range("BA9").Value = 70
$ws.Range("BB9").Value = 200

# Row 19: 24 cell(s) updated
$ws.Range("G19").Value = 1.6
$ws.Range("H19").Value = 4.5
$ws.Range("I19").Value = 4.75
$ws.Range("J19").Value = 2.05
$ws.Range("K19").Value = 2.6
$ws.Range("L19").Value = 4.75
$ws.Range("M19").Value = 1.02
$ws.Range("N19").Value = 21
$ws.Range("Q19").Value = 1.44
$ws.Range("R19").Value = 2.7
$ws.Range("U19").Value = 1.5
$ws.Range("V19").Value = 2.5
$ws.Range("AC19").Value = 21
$ws.Range("AD19").Value = 9.5
$ws.Range("AE19").Value = 13
$ws.Range("AF19").Value = 34
$ws.Range("AH19").Value = 19
$ws.Range("AI19").Value = 29
$ws.Range("AO19").Value = 8
$ws.Range("AP19").Value = 15
$ws.Range("AU19").Value = 7
$ws.Range("AW19").Value = 7
$ws.Range("AZ19").Value = 67
$ws.Range("BA19").Value = 67

# Row 20: 17 cell(s) updated
$ws.Range("G20").Value = 1.57
$ws.Range("K20").Value = 2.5
$ws.Range("Q20").Value = 1.6
$ws.Range("R20").Value = 2.3
$ws.Range("S20").Value = 1.29
$ws.Range("T20").Value = 3.5
$ws.Range("U20").Value = 1.67
$ws.Range("V20").Value = 2.1
$ws.Range("W20").Value = 9
$ws.Range("AD20").Value = 8.5
$ws.Range("AJ20").Value = 15
$ws.Range("AM20").Value = 34
$ws.Range("AT20").Value = 3.5
$ws.Range("AU20").Value = 7.5
$ws.Range("AV20").Value = 41
$ws.Range("BA20").Value = 81
$ws.Range("BC20").Value = 451

# Row 21: 28 cell(s) updated
$ws.Range("G21").Value = 1.91
$ws.Range("H21").Value = 3.9
$ws.Range("I21").Value = 3.6
$ws.Range("J21").Value = 2.4
$ws.Range("K21").Value = 2.5
$ws.Range("L21").Value = 3.75
$ws.Range("S21").Value = 1.25
$ws.Range("T21").Value = 3.75
$ws.Range("U21").Value = 1.5
$ws.Range("V21").Value = 2.5
$ws.Range("W21").Value = 11
$ws.Range("X21").Value = 12
$ws.Range("Y21").Value = 9
$ws.Range("AC21").Value = 19
$ws.Range("AE21").Value = 12
$ws.Range("AF21").Value = 34
$ws.Range("AG21").Value = 101
$ws.Range("AI21").Value = 21
$ws.Range("AJ21").Value = 12
$ws.Range("AL21").Value = 23
$ws.Range("AM21").Value = 26
$ws.Range("AO21").Value = 10
$ws.Range("AT21").Value = 3.75
$ws.Range("AX21").Value = 17
$ws.Range("AY21").Value = 21
$ws.Range("BA21").Value = 51
$ws.Range("BB21").Value = 101
$ws.Range("BC21").Value = 301

# Row 22: 15 cell(s) updated
$ws.Range("G22").Value = 1.83
$ws.Range("H22").Value = 4.1
$ws.Range("K22").Value = 2.5
$ws.Range("N22").Value = 19
$ws.Range("Q22").Value = 1.48
$ws.Range("R22").Value = 2.6
$ws.Range("S22").Value = 1.25
$ws.Range("T22").Value = 3.75
$ws.Range("AA22").Value = 13
$ws.Range("AB22").Value = 19
$ws.Range("AD22").Value = 8.5
$ws.Range("AG22").Value = 101
$ws.Range("AH22").Value = 17
$ws.Range("AT22").Value = 3.75
$ws.Range("AW22").Value = 6

# Row 23: 41 cell(s) updated
$ws.Range("G23").Value = 2.8
$ws.Range("H23").Value = 3.4
$ws.Range("I23").Value = 2.5
$ws.Range("J23").Value = 3.25
$ws.Range("K23").Value = 2.2
$ws.Range("L23").Value = 3.1
$ws.Range("M23").Value = 1.05
$ws.Range("N23").Value = 11
$ws.Range("O23").Value = 1.25
$ws.Range("P23").Value = 3.75
$ws.Range("Q23").Value = 1.88
$ws.Range("R23").Value = 1.98
$ws.Range("S23").Value = 1.36
$ws.Range("T23").Value = 3
$ws.Range("U23").Value = 1.67
$ws.Range("V23").Value = 2.1
$ws.Range("X23").Value = 15
$ws.Range("Y23").Value = 10
$ws.Range("Z23").Value = 29
$ws.Range("AA23").Value = 21
$ws.Range("AB23").Value = 29
$ws.Range("AC23").Value = 11
$ws.Range("AE23").Value = 13
$ws.Range("AH23").Value = 9.5
$ws.Range("AI23").Value = 13
$ws.Range("AK23").Value = 23
$ws.Range("AL23").Value = 19
$ws.Range("AN23").Value = 4.75
$ws.Range("AO23").Value = 15
$ws.Range("AP23").Value = 23
$ws.Range("AQ23").Value = 51
$ws.Range("AR23").Value = 67
$ws.Range("AS23").Value = 151
$ws.Range("AT23").Value = 3
$ws.Range("AV23").Value = 51
$ws.Range("AW23").Value = 4.5
$ws.Range("AX23").Value = 13
$ws.Range("AZ23").Value = 41
$ws.Range("BB23").Value = 151
$ws.Range("BC23").Value = 501
$ws.Range("BD23").Value = 126

# Row 24: 33 cell(s) updated
$ws.Range("G24").Value = 2.05
$ws.Range("H24").Value = 4.1
$ws.Range("J24").Value = 2.5
$ws.Range("K24").Value = 2.63
$ws.Range("M24").Value = 1.01
$ws.Range("N24").Value = 23
$ws.Range("O24").Value = 1.08
$ws.Range("P24").Value = 8
$ws.Range("Q24").Value = 1.33
$ws.Range("R24").Value = 3.4
$ws.Range("S24").Value = 1.2
$ws.Range("T24").Value = 4.33
$ws.Range("U24").Value = 1.33
$ws.Range("V24").Value = 3.25
$ws.Range("W24").Value = 17
$ws.Range("X24").Value = 17
$ws.Range("Y24").Value = 10
$ws.Range("AA24").Value = 13
$ws.Range("AC24").Value = 29
$ws.Range("AD24").Value = 9.5
$ws.Range("AF24").Value = 23
$ws.Range("AG24").Value = 67
$ws.Range("AH24").Value = 21
$ws.Range("AI24").Value = 23
$ws.Range("AN24").Value = 5
$ws.Range("AO24").Value = 10
$ws.Range("AP24").Value = 13
$ws.Range("AQ24").Value = 29
$ws.Range("AT24").Value = 4.33
$ws.Range("AW24").Value = 6
$ws.Range("BA24").Value = 41
$ws.Range("BB24").Value = 67
$ws.Range("BC24").Value = 151

# Row 25: 30 cell(s) updated
$ws.Range("H25").Value = 4.2
$ws.Range("J25").Value = 2.05
$ws.Range("K25").Value = 2.5
$ws.Range("L25").Value = 5.5
$ws.Range("M25").Value = 1.03
$ws.Range("N25").Value = 15
$ws.Range("O25").Value = 1.17
$ws.Range("P25").Value = 5
$ws.Range("Q25").Value = 1.57
$ws.Range("R25").Value = 2.35
$ws.Range("S25").Value = 1.29
$ws.Range("T25").Value = 3.5
$ws.Range("U25").Value = 1.67
$ws.Range("V25").Value = 2.1
$ws.Range("W25").Value = 9
$ws.Range("X25").Value = 8.5
$ws.Range("AB25").Value = 21
$ws.Range("AC25").Value = 15
$ws.Range("AD25").Value = 8
$ws.Range("AE25").Value = 15
$ws.Range("AG25").Value = 151
$ws.Range("AK25").Value = 67
$ws.Range("AM25").Value = 41
$ws.Range("AR25").Value = 41
$ws.Range("AS25").Value = 101
$ws.Range("AT25").Value = 3.5
$ws.Range("AU25").Value = 8
$ws.Range("AY25").Value = 29
$ws.Range("BA25").Value = 101
$ws.Range("BC25").Value = 451

# Row 26: 9 cell(s) updated
$ws.Range("G26").Value = 2.15
$ws.Range("I26").Value = 3.2
$ws.Range("O26").Value = 1.22
$ws.Range("P26").Value = 4
$ws.Range("Q26").Value = 1.73
$ws.Range("R26").Value = 2.08
$ws.Range("AO26").Value = 11
$ws.Range("AP26").Value = 19
$ws.Range("AW26").Value = 5.5

# Row 27: 11 cell(s) updated
$ws.Range("G27").Value = 2.38
$ws.Range("I27").Value = 2.63
$ws.Range("J27").Value = 2.75
$ws.Range("L27").Value = 2.88
$ws.Range("Y27").Value = 11
$ws.Range("AJ27").Value = 11
$ws.Range("AK27").Value = 29
$ws.Range("AL27").Value = 17
$ws.Range("AM27").Value = 17
$ws.Range("AO27").Value = 12
$ws.Range("AP27").Value = 15

# Row 28: 9 cell(s) updated
$ws.Range("G28").Value = 2
$ws.Range("I28").Value = 3.4
$ws.Range("J28").Value = 2.5
$ws.Range("L28").Value = 3.5
$ws.Range("N28").Value = 19
$ws.Range("U28").Value = 1.4
$ws.Range("V28").Value = 2.75
$ws.Range("AD28").Value = 8
$ws.Range("AI28").Value = 21

# Row 29: 12 cell(s) updated
$ws.Range("G29").Value = 1.42
$ws.Range("I29").Value = 5.75
$ws.Range("J29").Value = 1.83
$ws.Range("L29").Value = 5.5
$ws.Range("U29").Value = 1.5
$ws.Range("V29").Value = 2.5
$ws.Range("Z29").Value = 12
$ws.Range("AE29").Value = 15
$ws.Range("AG29").Value = 101
$ws.Range("AQ29").Value = 17
$ws.Range("AW29").Value = 8.5
$ws.Range("AZ29").Value = 81

# Row 32: 2 cell(s) updated
$ws.Range("Q32").Value = 1.88
$ws.Range("R32").Value = 1.98
